# ---------------------------------------------------------------------------
# Edit summary (see commit message / xml diff):
#   1. The table on the "B1- TYPES OF FINANCIAL DOCUMENTS" slide gets a new
#      table style applied: {AC74B768-E4C3-4D0E-BE17-E007ABBB1487} ->
#      {14D5D6D4-8559-4304-BC88-1CA2623DFC93}.
#   2. The deck's two themes swap places: the theme backing the slide master
#      (currently "Integral" / Red Violet) becomes the plain "Office" colour
#      palette, while the other theme (plain Office, used by the notes
#      master) becomes the Red Violet / Integral palette.
#      Table-styles and theme internals (fonts, format scheme, name
#      attributes) aren't exposed as settable properties anywhere in the
#      PowerPoint object model, but the 12 theme colours are (via
#      ThemeColorScheme), and this is the part of the swap that is visible
#      to slides, so it is applied color-by-color below.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Retarget the table style on the slide that has one.
# ---------------------------------------------------------------------------
$oldStyleId = "{AC74B768-E4C3-4D0E-BE17-E007ABBB1487}"
$newStyleId = "{14D5D6D4-8559-4304-BC88-1CA2623DFC93}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Swap the two theme colour palettes.
#
# Helper: convert a "RRGGBB" hex string into the little-endian integer that
# PowerPoint's RGB()/ThemeColor.RGB expect (R + G*256 + B*65536).
# ---------------------------------------------------------------------------
function ConvertTo-VbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order of ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officePalette = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47",
    "0563C1", "954F72"
)

# The slide master's theme currently holds the Red Violet / "Integral"
# palette; repaint it with the plain Office palette (mirrors the target
# theme2.xml content). Going through slide 1 reaches the shared master
# theme, so every slide/layout picks up the change.
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-VbaRgb($officePalette[$i - 1])
}
